# Generate Report for Archive
# 1. Update the "Ready for handoff" status text to "In Translation" everywhere it occurs
#    (Overview!E2:F2, zh-cn!C2, de-de!C2 all share this string).
# 2. Re-fit the Status / language columns that held that text: they shrink from the
#    "Ready for handoff"-sized width down to the narrower "In Translation" width.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $v = $cell.Value2
        # Guard with -is [string]: some cells hold native Booleans (e.g. "True"/"False"
        # status flags), and comparing a Boolean to a non-empty string with -eq coerces
        # the string to Boolean (true), producing a false-positive match.
        if (($v -is [string]) -and ($v -eq "Ready for handoff")) {
            $cell.Value2 = "In Translation"
        }
    }
}

# Target column width (character units) that corresponds to the new, narrower
# stored sheet width of ~13.41 once re-fit for "In Translation".
$newColumnWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $newColumnWidth
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = $newColumnWidth
